$wb = $excel.ActiveWorkbook

# VT-SaleVoid-DualCF-Generic
$ws = $wb.Worksheets.Item("VT-SaleVoid-DualCF-Generic")
$ws.Range("B2").Value = "Tue Aug 12 02:44:33 IST 2025"
$ws.Range("B3").Value = "Tue Aug 12 02:45:13 IST 2025"
$ws.Range("B4").Value = "Tue Aug 12 02:45:51 IST 2025"
$ws.Range("B5").Value = "Tue Aug 12 02:46:33 IST 2025"

# VT-SaleVoid-NoCF-Generic
$ws = $wb.Worksheets.Item("VT-SaleVoid-NoCF-Generic")
$ws.Range("D2").Value = "Tue Aug 12 02:47:15 IST 2025"
$ws.Range("D3").Value = "Tue Aug 12 02:47:56 IST 2025"
$ws.Range("D4").Value = "Tue Aug 12 02:48:43 IST 2025"
$ws.Range("D5").Value = "Tue Aug 12 02:49:27 IST 2025"

# VT-SaleVoid-SingleCF-Generic
$ws = $wb.Worksheets.Item("VT-SaleVoid-SingleCF-Generic")
$ws.Range("B2").Value = "Tue Aug 12 02:50:06 IST 2025"
$ws.Range("B3").Value = "Tue Aug 12 02:50:48 IST 2025"
$ws.Range("B4").Value = "Tue Aug 12 02:51:26 IST 2025"
$ws.Range("B5").Value = "Tue Aug 12 02:52:11 IST 2025"

# VT-SaleCredit-DualCF-Generic
$ws = $wb.Worksheets.Item("VT-SaleCredit-DualCF-Generic")
$ws.Range("B2").Value = "Tue Aug 12 02:36:25 IST 2025"
$ws.Range("B3").Value = "Tue Aug 12 02:37:05 IST 2025"
$ws.Range("B4").Value = "Tue Aug 12 02:37:54 IST 2025"
$ws.Range("B5").Value = "Tue Aug 12 02:38:35 IST 2025"

# VT-SaleCredit-NoCF-Generic
$ws = $wb.Worksheets.Item("VT-SaleCredit-NoCF-Generic")
$ws.Range("B2").Value = "Tue Aug 12 02:39:20 IST 2025"
$ws.Range("B3").Value = "Tue Aug 12 02:39:56 IST 2025"
$ws.Range("B4").Value = "Tue Aug 12 02:40:34 IST 2025"
$ws.Range("B5").Value = "Tue Aug 12 02:41:14 IST 2025"

# VT-SaleCredit-SingleCF-Generic
$ws = $wb.Worksheets.Item("VT-SaleCredit-SingleCF-Generic")
$ws.Range("B2").Value = "Tue Aug 12 02:41:52 IST 2025"
$ws.Range("B3").Value = "Tue Aug 12 02:42:33 IST 2025"
$ws.Range("B4").Value = "Tue Aug 12 02:43:10 IST 2025"
$ws.Range("B5").Value = "Tue Aug 12 02:43:54 IST 2025"

# VT-AuthCapCredit-Generic
$ws = $wb.Worksheets.Item("VT-AuthCapCredit-Generic")
$ws.Range("D2").Value = "Tue Aug 12 01:54:47 IST 2025"
$ws.Range("D3").Value = "Tue Aug 12 01:55:49 IST 2025"
$ws.Range("D4").Value = "Tue Aug 12 01:56:49 IST 2025"
$ws.Range("D5").Value = "Tue Aug 12 01:57:54 IST 2025"
$ws.Range("D6").Value = "Tue Aug 12 01:58:57 IST 2025"
$ws.Range("D7").Value = "Tue Aug 12 01:59:59 IST 2025"
$ws.Range("C5").Value = "Pass"

# VT-AuthCapVoid-Generic
$ws = $wb.Worksheets.Item("VT-AuthCapVoid-Generic")
$ws.Range("D2").Value = "Tue Aug 12 02:01:09 IST 2025"
$ws.Range("D3").Value = "Tue Aug 12 02:02:04 IST 2025"
$ws.Range("D4").Value = "Tue Aug 12 02:03:00 IST 2025"
$ws.Range("D5").Value = "Tue Aug 12 02:04:00 IST 2025"
$ws.Range("D6").Value = "Tue Aug 12 02:04:56 IST 2025"
$ws.Range("D7").Value = "Tue Aug 12 02:05:59 IST 2025"

# VT-ManualAuthCapture-Generic
$ws = $wb.Worksheets.Item("VT-ManualAuthCapture-Generic")
$ws.Range("B2").Value = "Tue Aug 12 02:32:14 IST 2025"
$ws.Range("B3").Value = "Tue Aug 12 02:32:59 IST 2025"
$ws.Range("B4").Value = "Tue Aug 12 02:33:38 IST 2025"
$ws.Range("B5").Value = "Tue Aug 12 02:34:20 IST 2025"
$ws.Range("B6").Value = "Tue Aug 12 02:34:58 IST 2025"
$ws.Range("B7").Value = "Tue Aug 12 02:35:42 IST 2025"
